$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Depth")
$tbl = $ws.ListObjects.Item("Table1")
$row = $tbl.ListRows.Add()
$ws.Range("A4").Value = "Pseudo 3D 2"
$ws.Range("B4").Value = 12
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 29.212765957446809
$ws.Range("E4").Value = 91.424177669598151
$ws.Range("F4").Value = 0.77645902037665526
$ws.Range("G4").Value = 0.7348843373594699
$ws.Range("H4").Value = 8
$ws.Range("I4").Value = "['NLD-AMPH-0017-LAD_frame452', 'NLD-AMPH-0063_frame297', 'NLD-AMPH-0063_frame440', 'NLD-ISALA-0084_frame59', 'NLD-ISALA-0084_frame360', 'NLD-RADB-0024_frame320', 'NLD-RADB-0024_frame520', 'NLD-RADB-0078_frame280']"
